# A new daily price record was inserted above the existing row 330 ("Fruta,
# Feria Lagunitas de Puerto Montt - Piña"), pushing the former rows 330:366
# down to 331:367 (the sheet's used range grows from A1:T366 to A1:T367).
#
# The new row re-uses the market/product context (Mercado, Región, Codreg,
# Tipo, Producto*, Categoría*, Variedad, Origen) that the old row 330 held,
# but carries its own date/quality/volume/price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the context columns of the current row 330 before the rows below
# shift down, so the freshly inserted row can reuse them.
$mercadoId = $ws.Range("A330").Value()
$mercado   = $ws.Range("B330").Value()
$region    = $ws.Range("C330").Value()
$codreg    = $ws.Range("E330").Value()
$tipo      = $ws.Range("F330").Value()
$prodId    = $ws.Range("G330").Value()
$producto  = $ws.Range("H330").Value()
$catId     = $ws.Range("I330").Value()
$categoria = $ws.Range("J330").Value()
$variedad  = $ws.Range("K330").Value()
$origen    = $ws.Range("R330").Value()

# Insert a new blank row above the current row 330; this shifts rows
# 330:366 down to 331:367 and grows the sheet dimension to A1:T367.
$ws.Rows.Item(330).Insert()

# Populate the newly inserted row 330.
$ws.Range("A330").Value = $mercadoId
$ws.Range("B330").Value = $mercado
$ws.Range("C330").Value = $region
$ws.Range("D330").Value = 44946
$ws.Range("E330").Value = $codreg
$ws.Range("F330").Value = $tipo
$ws.Range("G330").Value = $prodId
$ws.Range("H330").Value = $producto
$ws.Range("I330").Value = $catId
$ws.Range("J330").Value = $categoria
$ws.Range("K330").Value = $variedad
$ws.Range("L330").Value = "Segunda"
$ws.Range("M330").Value = 200
$ws.Range("N330").Value = 11000
$ws.Range("O330").Value = 12000
$ws.Range("P330").Value = 11500
$ws.Range("Q330").Value = "$/caja 7 unidades"
$ws.Range("R330").Value = $origen
$ws.Range("S330").Value = 1643
$ws.Range("T330").Value = 7
